$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @"
"严重性指该问题产生后果的严重程度。
满足下述任一维度，该问题即被判定为”严重性等级高“：
维度1-合规性风险，包括如“内部控制制度、是否符合相关制度、是否履行必要的审议程序、全面梳理、管理层是否勤勉尽责、是否损害投资者利益”等词
维度2-经营性风险，包括如“持续经营能力、终止上市风险、是否影响日常经营业务开展、退市风险”等词
维度3-重大偏离风险，监管机构对问题的描述存在“重大不确定性、远高于、远低于、大幅增长、大幅下降“等这类副词"
问询问题内容既不涉及合规性风险，也不涉及经营性风险和重大偏离风险，则该问题被判定为”严重性等级低“。
**返回结果请从["严重性等级高", "严重性等级低"]中枚举**
"@

$cell = $ws.Range("A2")
$cell.Value = $newText
$cell.WrapText = $true
$ws.Rows(2).RowHeight = 409.5

$cell.Select()
